$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C got a bit narrower (14.7109375 -> 13.7109375 width units) ---
$ws.Columns.Item(3).ColumnWidth = 12.833333333333334

# --- Re-apply the text format to the header row + Cell_class column so the ---
# --- redundant/unused cell-format records collapse onto the first text style ---
$ws.Range("A1:A10").NumberFormat = "@"
$ws.Range("B1:E1").NumberFormat = "@"

# --- Refreshed analysis numbers (re-run stats) ---
$ws.Range("B2").Value = 2.0890302066772657
$ws.Range("C2").Value = 0.018107054698453657
$ws.Range("D2").Value = 6.5748855325584756
$ws.Range("E2").Value = [double]"3.4351309756340717e-28"

$ws.Range("B3").Value = 0.43822666896670692
$ws.Range("C3").Value = 1.8018869873266381

$ws.Range("B5").Value = 1.1363900693357192
$ws.Range("C5").Value = 1.359762246723514

$ws.Range("B7").Value = 0.19306885544915642
$ws.Range("C7").Value = 1.7913027939052881

$ws.Range("B8").Value = 0.67925133689839567
$ws.Range("C8").Value = 1.7413094034507368

$ws.Range("B9").Value = 3.3877018817602607
$ws.Range("C9").Value = [double]"0.0065936140795095841"

$ws.Range("B10").Value = 5.9575532538596834
$ws.Range("C10").Value = [double]"1.1202021208138104e-11"
